# Append Week 17 NFL scores (rows 18-32) below the existing Week 16 data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Week, Date, Visitor, Visitor_pts, Home, Home_pts
$data = @(
    @(17, "12/25/2025", "Dallas Cowboys",        "30", "Washington Commanders", "23"),
    @(17, "12/25/2025", "Detroit Lions",          "10", "Minnesota Vikings",     "23"),
    @(17, "12/25/2025", "Denver Broncos",         "20", "Kansas City Chiefs",    "13"),
    @(17, "12/27/2025", "Houston Texans",         "20", "Los Angeles Chargers",  "16"),
    @(17, "12/27/2025", "Baltimore Ravens",       "41", "Green Bay Packers",     "24"),
    @(17, "12/28/2025", "Seattle Seahawks",       "27", "Carolina Panthers",     "10"),
    @(17, "12/28/2025", "Arizona Cardinals",      "14", "Cincinnati Bengals",    "37"),
    @(17, "12/28/2025", "New Orleans Saints",     "34", "Tennessee Titans",      "26"),
    @(17, "12/28/2025", "Jacksonville Jaguars",   "23", "Indianapolis Colts",    "17"),
    @(17, "12/28/2025", "Tampa Bay Buccaneers",   "17", "Miami Dolphins",        "20"),
    @(17, "12/28/2025", "New England Patriots",   "42", "New York Jets",         "10"),
    @(17, "12/28/2025", "Pittsburgh Steelers",     "6", "Cleveland Browns",      "13"),
    @(17, "12/28/2025", "New York Giants",        "34", "Las Vegas Raiders",     "10"),
    @(17, "12/28/2025", "Philadelphia Eagles",    "13", "Buffalo Bills",         "12"),
    @(17, "12/28/2025", "Chicago Bears",          "38", "San Francisco 49ers",   "42")
)

$startRow = 18
$endRow = $startRow + $data.Count - 1

# The Date (B) and points (D, F) columns must land as literal text, exactly
# like the rest of the sheet, rather than Excel auto-coercing them into a
# date serial / number. Temporarily force Text format before writing, then
# drop the formatting again so the new cells keep the sheet's default style.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"
$ws.Range("D$startRow`:D$endRow").NumberFormat = "@"
$ws.Range("F$startRow`:F$endRow").NumberFormat = "@"

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$ws.Range("B$startRow`:B$endRow").ClearFormats()
$ws.Range("D$startRow`:D$endRow").ClearFormats()
$ws.Range("F$startRow`:F$endRow").ClearFormats()
